# "Added MPI 1 process results" - append a new row (row 3) of benchmark
# results to the MPI sheet, mirroring the existing row 2 layout/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPI")
$ws.Activate()

# nodes / cores-per-node / matrix dim. (plain numbers, no special format)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Formula = "=A3*B3"
$ws.Range("D3").Value = 5000

# Timing/overhead columns use the same "0.00" number format as row 2.
$ws.Range("E3:M3").NumberFormat = "0.00"

$ws.Range("E3").Value = 2.245485
$ws.Range("F3").Value = 790.76326800000004
$ws.Range("G3").Value = 19.051117000000001
$ws.Range("H3").Value = 2.3856459999999999
$ws.Range("I3").Value = 1166.8051579999999
$ws.Range("J3").Value = 233.70181199999999
$ws.Range("K3").Formula = "=H3/E3"
$ws.Range("L3").Formula = "=I3/F3"
$ws.Range("M3").Formula = "=J3/G3"

# Matches the committed selection state: cursor on the first cell of the
# newly added row.
$ws.Range("A3").Select()
